# T31R09_data.xlsx restructuring:
#  - "1881"  data tab  -> renamed to "1880Survey"   (tab reads 1880 for consistency
#                          with the R data-processing script even though the survey
#                          was actually conducted in 1881)
#  - "1940"  data tab  -> renamed to "1940Survey"   (parallel naming convention)
#  - new "1880Metadata" tab appended, documenting the 1880/1881 naming quirk
#  - "1940Survey" becomes the active/selected tab when the workbook is saved

$wb = $excel.ActiveWorkbook

# --- rename the two data tabs to the new parallel "<year>Survey" convention ---
$ws1881 = $wb.Worksheets.Item("1881")
$ws1881.Name = "1880Survey"
$ws1881.Activate()
$ws1881.Range("A1").Select()

$ws1940 = $wb.Worksheets.Item("1940")
$ws1940.Name = "1940Survey"

# --- append a new metadata sheet explaining the 1880/1881 naming ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$metaSheet.Name = "1880Metadata"
$metaSheet.Range("A1").Value = "Actually surveyed in 1881; tab reads 1880 for consistency in the data processing R script"

# --- make the 1940Survey tab the active tab, matching the saved workbook view ---
$ws1940.Activate()
$ws1940.Range("A1").Select()

Write-Output ($wb.Worksheets | ForEach-Object { $_.Name })
